$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for WY2024->WY2025 "Delta" data (becomes row 7) ---
# This shifts the existing "Benicia" rows (formerly rows 7-12) down to rows 8-13.
$ws.Rows.Item(7).Insert()

# Copy formatting down from the row above (row 6) so the new row matches the
# existing "Delta" block styling (font/border per column).
$ws.Range("A6:H6").Copy() | Out-Null
$ws.Range("A7:H7").PasteSpecial(-4122) | Out-Null

$ws.Range("A7").Value2 = "Delta"
$ws.Range("B7").Value2 = 2024
$ws.Range("C7").Value2 = 2025
$ws.Range("D7").Value2 = 64.1
$ws.Range("E7").Value2 = 2.8
$ws.Range("F7").Value2 = 58.5
$ws.Range("G7").Value2 = 69.3
$ws.Range("H7").ClearContents() | Out-Null
$ws.Rows.Item(7).RowHeight = 15.75

# --- Extend the WY shared formula down through the new last "Benicia" row (13) ---
$ws.Range("B13").Formula = "=C13-1"

# --- Append a new "Benicia" row for WY2024->WY2025 (row 14) ---
$ws.Range("A13:H13").Copy() | Out-Null
$ws.Range("A14:H14").PasteSpecial(-4122) | Out-Null

$ws.Range("A14").Value2 = "Benicia"
$ws.Range("B14").Value2 = 2024
$ws.Range("C14").Value2 = 2025
$ws.Range("D14").Value2 = 26.7
$ws.Range("E14").Value2 = 1.5
$ws.Range("F14").Value2 = 23.8
$ws.Range("G14").Value2 = 29.7
$ws.Range("H14").Value2 = 93.8
$ws.Rows.Item(14).RowHeight = 15.75

$ws.Range("H8").Select() | Out-Null
